# ---------------------------------------------------------------------------
# Commit: "feat: add 2022-Q3 data"
#
# 1) Insert a new "2022-Q3" row at the top of the "总计" (summary) sheet,
#    shifting the existing quarters down and renumbering the index column.
# 2) Insert a brand-new worksheet named "2022-Q3" right after "总计" holding
#    the per-fund holding breakdown for that quarter (mirrors the layout of
#    the existing quarter sheets).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Update the "总计" summary sheet -------------------------------------
$summary = $wb.Worksheets.Item(1)

# Make room for the new quarter at row 2 (push 2022-Q2/2022-Q1/2021-Q1/2020-Q4
# down by one row each), copying formatting from the row being displaced so
# the inserted row keeps the same style as the rest of the table.
$summary.Rows.Item(2).Insert()
$summary.Range("A3:D3").Copy()
$summary.Range("A2:D2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 49
$summary.Range("D2").Value = 101.85

# Renumber the index column for the rows that shifted down.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4

# --- 2. Add the new "2022-Q3" worksheet -------------------------------------
# Insert it right after "总计", i.e. before "2022-Q2".
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

# Header row (matches the layout used by the other quarter sheets).
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$q3Data = @(
    @("0", "513050", "易方达中证海外中国互联网50（QDII）ETF", "317.38", "98.16", "7.52", "23.8670", "4"),
    @("1", "513330", "华夏恒生互联网科技业ETF（QDII）", "204.59", "95.71", "11.13", "22.7709", "4"),
    @("2", "110011", "易方达优质精选混合（QDII）", "173.81", "94.25", "7.86", "13.6615", "6"),
    @("3", "513180", "华夏恒生科技ETF（QDII）", "103.92", "94.41", "7.41", "7.7005", "4"),
    @("4", "164906", "交银施罗德中证海外中国互联网指数（QDII-LOF）", "97.68", "91.19", "6.83", "6.6715", "4"),
    @("5", "159920", "华夏恒生ETF（QDII）", "147.66", "95.67", "3.00", "4.4298", "7"),
    @("6", "510900", "易方达恒生国企（QDII-ETF）", "81.42", "95.82", "5.13", "4.1768", "5"),
    @("7", "118001", "易方达亚洲精选股票（QDII）", "46.17", "94.52", "8.25", "3.8090", "1"),
    @("8", "159605", "广发中证海外中国互联网30（QDII-ETF）", "33.01", "99.15", "10.12", "3.3406", "4"),
    @("9", "513010", "易方达恒生科技ETF（QDII）", "25.49", "94.56", "7.45", "1.8990", "4"),
    @("10", "070012", "嘉实海外中国混合（QDII）", "16.90", "82.25", "5.60", "0.9464", "6"),
    @("11", "159740", "大成恒生科技ETF（QDII）", "9.95", "93.00", "7.33", "0.7293", "4"),
    @("12", "159607", "嘉实中证海外中国互联网30ETF（QDII）", "6.55", "98.57", "10.05", "0.6583", "4"),
    @("13", "000988", "嘉实全球互联网股票-人民币（QDII）", "10.62", "88.35", "4.91", "0.5214", "6"),
    @("14", "000989", "嘉实全球互联网股票-美元现汇（QDII）", "10.62", "88.35", "4.91", "0.5214", "6"),
    @("15", "000990", "嘉实全球互联网股票-美元现钞（QDII）", "10.62", "88.35", "4.91", "0.5214", "6"),
    @("16", "012805", "广发恒生科技指数（QDII）C", "7.98", "80.85", "6.37", "0.5083", "4"),
    @("17", "513580", "华安恒生科技ETF（QDII）", "6.44", "95.77", "7.47", "0.4811", "4"),
    @("18", "159742", "博时恒生科技ETF（QDII）", "6.13", "95.25", "7.51", "0.4604", "4"),
    @("19", "010671", "景顺长城大中华混合（QDII）美元A", "9.42", "70.56", "3.64", "0.3429", "8"),
    @("20", "262001", "景顺长城大中华混合（QDII）人民币A", "9.42", "70.56", "3.64", "0.3429", "8"),
    @("21", "012208", "华夏港股前沿经济混合（QDII）A", "9.67", "89.48", "3.42", "0.3307", "9"),
    @("22", "159747", "南方中证香港科技ETF（QDII）", "2.92", "99.57", "10.95", "0.3197", "2"),
    @("23", "862001", "光大阳光香港精选混合（QDII）A 人民币", "3.15", "90.62", "7.62", "0.2400", "3"),
    @("24", "862011", "光大阳光香港精选混合（QDII）A 美元", "3.15", "90.62", "7.62", "0.2400", "3"),
    @("25", "862012", "光大阳光香港精选混合（QDII）C 人民币", "3.15", "90.62", "7.62", "0.2400", "3"),
    @("26", "159850", "华夏恒生中国企业ETF（QDII）", "4.84", "93.33", "4.86", "0.2352", "5"),
    @("27", "012379", "创金合信港股互联网3个月持有期混合（QDII）A", "2.81", "87.48", "8.08", "0.2270", "3"),
    @("28", "159741", "嘉实恒生科技ETF（QDII）", "2.85", "99.55", "7.85", "0.2237", "4"),
    @("29", "012804", "广发恒生科技指数（QDII）A", "3.18", "80.85", "6.37", "0.2026", "4"),
    @("30", "513890", "上投摩根恒生科技ETF（QDII）", "1.73", "95.47", "8.44", "0.1460", "1"),
    @("31", "160717", "嘉实恒生中国企业指数（QDII-LOF）", "2.42", "94.01", "5.03", "0.1217", "5"),
    @("32", "513220", "招商中证全球中国互联网ETF（QDII）", "1.05", "98.86", "11.13", "0.1169", "4"),
    @("33", "013127", "汇添富恒生科技指数（QDII）A", "1.49", "91.29", "7.20", "0.1073", "4"),
    @("34", "013128", "汇添富恒生科技指数（QDII）C", "1.45", "91.29", "7.20", "0.1044", "4"),
    @("35", "164705", "汇添富恒生指数（QDII-LOF）A", "2.74", "90.36", "3.34", "0.0915", "7"),
    @("36", "012380", "创金合信港股互联网3个月持有期混合（QDII）C", "0.96", "87.48", "8.08", "0.0776", "3"),
    @("37", "159750", "招商中证香港科技ETF（QDII）", "0.69", "98.87", "10.85", "0.0749", "2"),
    @("38", "513380", "广发恒生科技（QDII-ETF）", "1.08", "84.71", "6.81", "0.0735", "3"),
    @("39", "160125", "南方香港优选股票（QDII-LOF）", "2.01", "81.74", "3.12", "0.0627", "8"),
    @("40", "160644", "鹏华香港美国互联网股票（LOF）人民币", "1.23", "83.13", "4.40", "0.0541", "6"),
    @("41", "006792", "鹏华香港美国互联网股票（LOF）美元现汇", "1.23", "83.13", "4.40", "0.0541", "6"),
    @("42", "161229", "国投瑞银中国价值发现股票（QDII-LOF）", "1.33", "93.58", "3.62", "0.0481", "7"),
    @("43", "519601", "海富通中国海外精选混合（QDII）", "0.51", "73.52", "6.27", "0.0320", "1"),
    @("44", "160924", "大成恒生指数（QDII-LOF）", "0.90", "91.23", "2.87", "0.0258", "7"),
    @("45", "010789", "汇添富恒生指数（QDII-LOF）C", "0.43", "90.36", "3.34", "0.0144", "7"),
    @("46", "012209", "华夏港股前沿经济混合（QDII）C", "0.35", "89.48", "3.42", "0.0120", "9"),
    @("47", "378006", "上投摩根全球新兴市场混合（QDII）", "0.40", "87.48", "2.60", "0.0104", "7"),
    @("48", "519602", "海富通大中华精选混合（QDII）", "0.10", "87.37", "6.32", "0.0063", "2")
)

# Columns B-G hold text that looks numeric (fund codes, percentages, etc.) in
# every other quarter sheet, so force them to Text before writing so they
# aren't auto-converted to numbers.
$q3.Range("B2:G50").NumberFormat = "@"

for ($i = 0; $i -lt $q3Data.Count; $i++) {
    $row = $q3Data[$i]
    $r = $i + 2
    $q3.Cells.Item($r, 1).Value = [double]$row[0]
    $q3.Cells.Item($r, 2).Value = $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = $row[3]
    $q3.Cells.Item($r, 5).Value = $row[4]
    $q3.Cells.Item($r, 6).Value = $row[5]
    $q3.Cells.Item($r, 7).Value = $row[6]
    $q3.Cells.Item($r, 8).Value = [double]$row[7]
}
